$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: person responsible for "email html design/megvalósítás" changed from Barta to Géczy
$ws.Range("B7").Value = "Géczy"

# New row 9: elfelejtett jelszó backend (forgotten password backend) - Barta, 5 hours
$ws.Range("A9").Value = "elfelejtett jelszó backend"
$ws.Range("B9").Value = "Barta"
$ws.Range("C9").Value = 5

# New row 10: elfelejtett jelszó frontend (forgotten password frontend) - Barta, 1 hour
$ws.Range("A10").Value = "elfelejtett jelszó frontend"
$ws.Range("B10").Value = "Barta"
$ws.Range("C10").Value = 1

# Move selection to reflect where the user ended up after editing
$ws.Range("E14").Select()
